$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift milestone names up: Milestone 4 -> row2 becomes "Milestone 2",
# Milestone 5 -> row3 becomes "Milestone 4", and the old row 4 (Milestone 6) is removed.
$ws.Range("A2").Value = "Milestone 2"
$ws.Range("A3").Value = "Milestone 4"

# Remove the now-obsolete last row entirely (was Milestone 6), shifting the
# used range back down to A1:B3.
$ws.Rows("4:4").Delete()
